$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper cells (D column) - preserve literal formatted numeric-looking strings
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "37.505.02"
$ws.Range("E2").Value = "  +2.96%  "
Set-TextValue "D3" "2.071.30"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "235.06"
$ws.Range("E5").Value = "  +0.22%  "
Set-TextValue "D6" "0.619"
$ws.Range("E6").Value = "  +3.41%  "
Set-TextValue "D7" "58.32"
$ws.Range("E7").Value = "  +6.49%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.384"
$ws.Range("E9").Value = "  +3.69%  "
Set-TextValue "D10" "59.07"
$ws.Range("E10").Value = "  +1.21%  "
Set-TextValue "D11" "0.0762"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("E12").Value = "  +3.78%  "
Set-TextValue "D13" "2.372.41"
$ws.Range("E13").Value = "  +3.56%  "
Set-TextValue "D14" "14.54"
$ws.Range("E14").Value = "  +2.51%  "
Set-TextValue "D15" "21.13"
$ws.Range("E15").Value = "  +4.10%  "
Set-TextValue "D16" "0.779"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("E17").Value = "  +2.08%  "
Set-TextValue "D18" "2.074.55"
$ws.Range("E18").Value = "  +3.72%  "
Set-TextValue "D19" "37.674.48"
$ws.Range("E19").Value = "  +3.02%  "
Set-TextValue "D20" "6.21"
$ws.Range("E20").Value = "  +17.66%  "
Set-TextValue "D21" "70.20"
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("E26").Value = "  +1.06%  "
Set-TextValue "D27" "166.75"
$ws.Range("E27").Value = "  +2.35%  "
Set-TextValue "D28" "1.51"
$ws.Range("E28").Value = "  +9.23%  "
Set-TextValue "D29" "8.91"
$ws.Range("E29").Value = "  +2.91%  "
Set-TextValue "D30" "19.27"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("E31").Value = "  +1.58%  "
Set-TextValue "D32" "0.118"
$ws.Range("E32").Value = "  +1.85%  "
Set-TextValue "D33" "4.51"
$ws.Range("E33").Value = "  +3.83%  "
Set-TextValue "D34" "0.0623"
$ws.Range("E34").Value = "  +3.37%  "
Set-TextValue "D35" "2.57"
$ws.Range("E35").Value = "  +8.85%  "
Set-TextValue "D36" "4.56"
$ws.Range("E36").Value = "  +7.29%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "3.37"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D38" "1.00"
$ws.Range("E38").Value = "  -0.10%  "
Set-TextValue "D39" "1.78"
$ws.Range("E39").Value = "  +1.33%  "
Set-TextValue "D40" "5.86"
$ws.Range("E40").Value = "  +4.41%  "
Set-TextValue "D41" "4.62"
$ws.Range("E41").Value = "  +26.89%  "
Set-TextValue "D42" "2.96"
$ws.Range("E42").Value = "  -1.10%  "
Set-TextValue "D43" "0.0949"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D44" "1.18"
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D45" "1.456.53"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "95.75"
$ws.Range("E46").Value = "  +7.25%  "
$ws.Range("E47").Value = "  +4.74%  "
Set-TextValue "D48" "15.84"
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("E49").Value = "  +4.19%  "
Set-TextValue "D50" "7.28"
$ws.Range("E50").Value = "  +6.41%  "
$ws.Range("E51").Value = "  +1.81%  "
